$wb = $excel.ActiveWorkbook

# Update "展览" (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1070
$ws1.Range("F3").Value = 29

# Update "全部类型" (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1070
$ws4.Range("F3").Value = 29
